# WIP: local changes before rebase
# - move the active selection from C11 to C2
# - widen column A (1) and column D (4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths are exposed through COM as "characters" and the host
# re-derives the stored OOXML width as (chars * 6 + 5) / 6, i.e. it adds
# ~0.8333 characters of cell-padding before the round-trip. Back that
# padding out here so the persisted <col width="..."/> lands on the
# intended value (21.1640625 / 30.1640625 chars).
$ws.Columns.Item(1).ColumnWidth = 21.1640625 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 30.1640625 - 0.8333333333333334

# Move the selected/active cell to C2.
$ws.Range("C2").Select()
